$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21
$lastRow = 65

# --- 1. Rename header cells: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into a native Excel Table (ListObject) ---
# The header row already carries direct formatting (bold/fill/border). Adding a
# ListObject directly on top of that range would make Excel capture the existing
# header look as a one-off "dxf" (headerRowDxfId) in styles.xml. To keep styles.xml
# byte-for-byte identical, stash the header formatting away, clear it, create the
# table, and then restore the exact same formatting via copy/paste (which re-uses
# the pre-existing style record instead of minting a differential format record).
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

$headerRange.Copy($scratchRange)
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$table.Name = "Table1"

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)
$scratchRange.Clear()
